# Implements the "two types of payment methods" flow change:
#  - drop the now-unused second sheet (PURCHASEDRESS2)
#  - rename the two test cases to reflect cheque / bank-wire payment
#  - widen column A on the remaining sheet
#  - re-point the dress_name cell in row 3 to the "Summer" dress and pick
#    up the same (default) style already used by the rest of row 3

$wb = $excel.ActiveWorkbook

# 1. Remove the PURCHASEDRESS2 worksheet entirely.
$ws2 = $wb.Worksheets.Item("PURCHASEDRESS2")
$deleted = $ws2.Delete()

$ws = $wb.Worksheets.Item("PURCHASEDRESS")

# 2. Rename the two test cases: from generic "...FROM_SEARCH" to the two
#    concrete payment flows (cheque / bank wire).
$ws.Range("A2").Value = "TC01_CUSTOMER_PLACES_ORDER_BY_CHEQUE"
$ws.Range("A3").Value = "TC02_CUSTOMER_PLACES_ORDER_BY_BANK_WIRE"

# 3. Widen column A to fit the longer test case names.
$ws.Columns.Item(1).ColumnWidth = 50.14

# 4. Row 2's TestCaseName cell picks up row 3's formatting (same style
#    already used elsewhere in the sheet).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# 5. C3 switches from the Chiffon sample row to the Summer dress sample,
#    matching the formatting already used in column A of that row.
$ws.Range("C3").Value = "Printed Summer Dress"
$ws.Range("A3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$excel.CutCopyMode = $false
